# Convert MB_ESTOP_STATE to enum stop reason
#
# 1. Rename the original sheet to "Modbus Map" and add a new
#    "E-Stop Reasons" lookup sheet after it.
# 2. Update the MB_ESTOP_STATE row (row 54) on the Modbus Map sheet:
#      - Description -> "Status of emergency stop (first trip reason)"
#      - Units       -> "enumeration"
#      - Notes       -> "See 'E-Stop Reasons' sheet.  Cleared when MB_RESET_ESTOP activated"
# 3. Populate the new "E-Stop Reasons" sheet with the enumeration values.

$wb = $excel.ActiveWorkbook

# --- Sheets: rename + add -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Modbus Map"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "E-Stop Reasons"

# --- Modbus Map: MB_ESTOP_STATE row (row 54) ------------------------------
$ws1.Range("C54").Value = "Status of emergency stop (first trip reason)"
$ws1.Range("E54").Value = "enumeration"
$ws1.Range("G54").Value = "See 'E-Stop Reasons' sheet.  Cleared when MB_RESET_ESTOP activated"

# --- E-Stop Reasons sheet content ------------------------------------------
$ws2.Range("B2").Value = "E-Stop Reasons"
$ws2.Range("B2").Font.Bold = $true

$reasons = @(
    "NOT_ESTOPPED,",
    "ESTOP_REMOTE_COMMAND,",
    "ESTOP_CURRENT_LIMIT_INWARD,",
    "ESTOP_CURRENT_LIMIT_OUTWARD,",
    "ESTOP_BATT_OVERVOLTAGE,",
    "ESTOP_EXTENSION_LIMIT_INWARD,",
    "ESTOP_EXTENSION_LIMIT_OUTWARD,",
    "ESTOP_ENCODER_FAILURE,",
    "ESTOP_HEARTBEAT_TIMEOUT"
)

for ($i = 0; $i -lt $reasons.Length; $i++) {
    $row = 4 + $i
    $ws2.Cells.Item($row, 2).Value = $i
    $ws2.Cells.Item($row, 3).Value = $reasons[$i]
}

# --- Restore focus back on the main sheet, matching the commit's selection -
[void]$ws1.Activate()
[void]$ws1.Range("C54").Select()

Write-Output "MB_ESTOP_STATE enum conversion applied"
